# "Changes of 22nd June 2022"
# Shift the logged date/time for this job row from 17-Jun-2022 to 22-Jun-2022
# (date serial 44729 -> 44734; time-of-day 07:00/06:45 -> 22:00/21:45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start/creation date+time columns
$ws.Range("P2").Value = 44734
$ws.Range("S2").Value = 0.91666666666666663
$ws.Range("T2").Value = 0.90625

# Duplicate date/time columns further along the row (CA2/CB2)
$ws.Range("CA2").Value = 44734
$ws.Range("CB2").Value = 0.91666666666666663

# Move the scroll position / selection to the left portion of the sheet,
# with T2 as the active cell.
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("T2").Select()
